$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TilePath values (column I) for specific rows, using new shared strings
# Order matters for shared string table creation: tatamuviliage(48), smallforest(49), farm(50)
$ws.Range("I10").Value = "tatamuviliage"
$ws.Range("I16").Value = "smallforest"
$ws.Range("I13").Value = "farm"

# Update Level (column C) for row 10
$ws.Range("C10").Value = 1

# Update the selected cell on the sheet
$ws.Range("A11").Select()
